{"js": "// Restyle the built-in \"Date\" paragraph style:\n//   - paragraph alignment: center -> right\n//   - add bold character formatting\nconst dateStyle = context.document.getStyles().getByNameOrNullObject(\"Date\");\ndateStyle.load(\"nameLocal\");\nawait context.sync();\n\nif (dateStyle.isNullObject) {\n  throw new Error('Style \"Date\" was not found in this document.');\n}\n\ndateStyle.paragraphFormat.alignment = Word.Alignment.right;\ndateStyle.font.bold = true;\n\nawait context.sync();\n", "ps1": "# Restyle the built-in \"Date\" paragraph style:\n#   - paragraph alignment: center -> right\n#   - add bold character formatting\n$d = $word.ActiveDocument\n$dateStyle = $d.Styles(\"Date\")\n$dateStyle.ParagraphFormat.Alignment = [Microsoft.Office.Interop.Word.WdParagraphAlignment]::wdAlignParagraphRight\n$dateStyle.Font.Bold = $true\n"}
